# Updated cryptos list on Tue Jul 16 14:12:47 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns of the crypto table,
# and fixes the ordering of the RenderToken / EthereumClassic rows (33/34).
# Numeric-looking price strings are written with a leading apostrophe so
# Excel keeps them as literal text (matching the original inlineStr cells)
# instead of silently reinterpreting/rounding them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.243.27"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "3.385.34"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'567.42"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Value = "'155.13"
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.384.38"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").Value = "'0.566"
$ws.Range("E9").Value = "  +7.23%  "
$ws.Range("D10").Value = "'7.30"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "'0.122"
$ws.Range("E11").Value = "  +3.49%  "
$ws.Range("D12").Value = "'0.436"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "3.975.45"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").Value = "'0.0000189"
$ws.Range("E15").Value = "  +4.93%  "
$ws.Range("D16").Value = "'27.28"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "63.232.83"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "3.390.49"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "'6.28"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "'14.04"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "'379.77"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "'8.00"
$ws.Range("E22").Value = "  -4.84%  "
$ws.Range("D23").Value = "'0.994"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").Value = "'72.09"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").Value = "'0.531"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'0.0000117"
$ws.Range("E26").Value = "  +21.75%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'6.11"
$ws.Range("E30").Value = "  +9.64%  "
$ws.Range("E31").Value = "  +5.25%  "
$ws.Range("D32").Value = "'2.00"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'23.25"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'6.41"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("D36").Value = "'6.83"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").Value = "'159.55"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").Value = "'1.45"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").Value = "2.913.60"
$ws.Range("E39").Value = "  +3.54%  "
$ws.Range("D40").Value = "'0.0755"
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("D41").Value = "'1.82"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").Value = "'26.71"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("E43").Value = "  -5.01%  "
$ws.Range("D44").Value = "'41.86"
$ws.Range("E44").Value = "  +2.26%  "
$ws.Range("D45").Value = "'4.35"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").Value = "'0.753"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "'23.13"
$ws.Range("E47").Value = "  +5.30%  "
$ws.Range("D48").Value = "'1.05"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("D49").Value = "'2.17"
$ws.Range("E49").Value = "  +20.56%  "
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("E51").Value = "  +1.07%  "
